$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows are all dated 2025-06-28 (serial 45836). Copy A344s
# date-number-format down first, so the new date cells reuse the
# existing style instead of minting a new one.
$ws.Range("A344").Copy() | Out-Null
$ws.Range("A345:A351").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$newDate = 45836

# Row 345
$ws.Cells.Item(345, 1).Value = $newDate
$ws.Cells.Item(345, 2).Value = "Flowering"
$ws.Cells.Item(345, 3).Value = "Large"
$ws.Cells.Item(345, 4).Value = 69
$ws.Cells.Item(345, 5).Value = 82
$ws.Cells.Item(345, 6).Formula = "=ABS(D345-E345)"
$ws.Cells.Item(345, 7).Value = 5.57
$ws.Cells.Item(345, 8).Value = 0.6
$ws.Cells.Item(345, 9).Value = "No"
$ws.Cells.Item(345, 10).Value = 2
$ws.Cells.Item(345, 11).Value = "Dark"
$ws.Cells.Item(345, 12).Value = 8
$ws.Cells.Item(345, 13).Value = 0.94
$ws.Cells.Item(345, 14).Value = 71
$ws.Cells.Item(345, 15).Value = 30.03
$ws.Cells.Item(345, 16).Value = 6
$ws.Cells.Item(345, 17).Value = 0.78
$ws.Cells.Item(345, 18).Value = 5.3
$ws.Cells.Item(345, 19).Value = 47
$ws.Cells.Item(345, 20).Value = 0

# Row 346
$ws.Cells.Item(346, 1).Value = $newDate
$ws.Cells.Item(346, 2).Value = "Nonflowering"
$ws.Cells.Item(346, 3).Value = "Medium"
$ws.Cells.Item(346, 4).Value = 69
$ws.Cells.Item(346, 5).Value = 82
$ws.Cells.Item(346, 6).Formula = "=ABS(D346-E346)"
$ws.Cells.Item(346, 7).Value = 5.57
$ws.Cells.Item(346, 8).Value = 1.2
$ws.Cells.Item(346, 9).Value = "No"
$ws.Cells.Item(346, 10).Value = 3
$ws.Cells.Item(346, 11).Value = "Dark"
$ws.Cells.Item(346, 12).Value = 8
$ws.Cells.Item(346, 13).Value = 0.94
$ws.Cells.Item(346, 14).Value = 71
$ws.Cells.Item(346, 15).Value = 30.03
$ws.Cells.Item(346, 16).Value = 6
$ws.Cells.Item(346, 17).Value = 0.78
$ws.Cells.Item(346, 18).Value = 5.3
$ws.Cells.Item(346, 19).Value = 47
$ws.Cells.Item(346, 20).Value = 0

# Row 347
$ws.Cells.Item(347, 1).Value = $newDate
$ws.Cells.Item(347, 2).Value = "Nonflowering"
$ws.Cells.Item(347, 3).Value = "Small"
$ws.Cells.Item(347, 4).Value = 69
$ws.Cells.Item(347, 5).Value = 82
$ws.Cells.Item(347, 6).Formula = "=ABS(D347-E347)"
$ws.Cells.Item(347, 7).Value = 5.57
$ws.Cells.Item(347, 8).Value = 1.1
$ws.Cells.Item(347, 9).Value = "No"
$ws.Cells.Item(347, 10).Value = 3
$ws.Cells.Item(347, 11).Value = "Dark"
$ws.Cells.Item(347, 12).Value = 8
$ws.Cells.Item(347, 13).Value = 0.94
$ws.Cells.Item(347, 14).Value = 71
$ws.Cells.Item(347, 15).Value = 30.03
$ws.Cells.Item(347, 16).Value = 6
$ws.Cells.Item(347, 17).Value = 0.78
$ws.Cells.Item(347, 18).Value = 5.3
$ws.Cells.Item(347, 19).Value = 47
$ws.Cells.Item(347, 20).Value = 0

# Row 348
$ws.Cells.Item(348, 1).Value = $newDate
$ws.Cells.Item(348, 2).Value = "Nonflowering"
$ws.Cells.Item(348, 3).Value = "Medium"
$ws.Cells.Item(348, 4).Value = 69
$ws.Cells.Item(348, 5).Value = 82
$ws.Cells.Item(348, 6).Formula = "=ABS(D348-E348)"
$ws.Cells.Item(348, 7).Value = 5.57
$ws.Cells.Item(348, 8).Value = 1.6
$ws.Cells.Item(348, 9).Value = "No"
$ws.Cells.Item(348, 10).Value = 3
$ws.Cells.Item(348, 11).Value = "Dark"
$ws.Cells.Item(348, 12).Value = 8
$ws.Cells.Item(348, 13).Value = 0.94
$ws.Cells.Item(348, 14).Value = 71
$ws.Cells.Item(348, 15).Value = 30.03
$ws.Cells.Item(348, 16).Value = 6
$ws.Cells.Item(348, 17).Value = 0.78
$ws.Cells.Item(348, 18).Value = 5.3
$ws.Cells.Item(348, 19).Value = 47
$ws.Cells.Item(348, 20).Value = 0

# Row 349
$ws.Cells.Item(349, 1).Value = $newDate
$ws.Cells.Item(349, 2).Value = "Nonflowering"
$ws.Cells.Item(349, 3).Value = "Medium"
$ws.Cells.Item(349, 4).Value = 69
$ws.Cells.Item(349, 5).Value = 82
$ws.Cells.Item(349, 6).Formula = "=ABS(D349-E349)"
$ws.Cells.Item(349, 7).Value = 5.57
$ws.Cells.Item(349, 8).Value = 1.5
$ws.Cells.Item(349, 9).Value = "No"
$ws.Cells.Item(349, 10).Value = 3
$ws.Cells.Item(349, 11).Value = "Dark"
$ws.Cells.Item(349, 12).Value = 8
$ws.Cells.Item(349, 13).Value = 0.94
$ws.Cells.Item(349, 14).Value = 71
$ws.Cells.Item(349, 15).Value = 30.03
$ws.Cells.Item(349, 16).Value = 6
$ws.Cells.Item(349, 17).Value = 0.78
$ws.Cells.Item(349, 18).Value = 5.3
$ws.Cells.Item(349, 19).Value = 47
$ws.Cells.Item(349, 20).Value = 0

# Row 350
$ws.Cells.Item(350, 1).Value = $newDate
$ws.Cells.Item(350, 2).Value = "Nonflowering"
$ws.Cells.Item(350, 3).Value = "Large"
$ws.Cells.Item(350, 4).Value = 69
$ws.Cells.Item(350, 5).Value = 82
$ws.Cells.Item(350, 6).Formula = "=ABS(D350-E350)"
$ws.Cells.Item(350, 7).Value = 5.57
$ws.Cells.Item(350, 8).Value = 4
$ws.Cells.Item(350, 9).Value = "No"
$ws.Cells.Item(350, 10).Value = 4
$ws.Cells.Item(350, 11).Value = "Dark"
$ws.Cells.Item(350, 12).Value = 8
$ws.Cells.Item(350, 13).Value = 0.94
$ws.Cells.Item(350, 14).Value = 71
$ws.Cells.Item(350, 15).Value = 30.03
$ws.Cells.Item(350, 16).Value = 6
$ws.Cells.Item(350, 17).Value = 0.78
$ws.Cells.Item(350, 18).Value = 5.3
$ws.Cells.Item(350, 19).Value = 47
$ws.Cells.Item(350, 20).Value = 0

# Row 351
$ws.Cells.Item(351, 1).Value = $newDate
$ws.Cells.Item(351, 2).Value = "Tree"
$ws.Cells.Item(351, 3).Value = "Medium"
$ws.Cells.Item(351, 4).Value = 69
$ws.Cells.Item(351, 5).Value = 82
$ws.Cells.Item(351, 6).Formula = "=ABS(D351-E351)"
$ws.Cells.Item(351, 7).Value = 5.57
$ws.Cells.Item(351, 8).Value = 10
$ws.Cells.Item(351, 9).Value = "No"
$ws.Cells.Item(351, 10).Value = 1
$ws.Cells.Item(351, 11).Value = "Dark"
$ws.Cells.Item(351, 12).Value = 8
$ws.Cells.Item(351, 13).Value = 0.94
$ws.Cells.Item(351, 14).Value = 71
$ws.Cells.Item(351, 15).Value = 30.03
$ws.Cells.Item(351, 16).Value = 6
$ws.Cells.Item(351, 17).Value = 0.78
$ws.Cells.Item(351, 18).Value = 5.3
$ws.Cells.Item(351, 19).Value = 47
$ws.Cells.Item(351, 20).Value = 0

# Match the author's final selection/view state
$ws.Range("AD11").Select()
